$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Cells.Item(8,1).Value = "Volume 30   Number  9"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  2/27/2023  Through  3/5/2023"

# --- Cells changing from numeric to special text markers ("0" / "***.*") ---
# Donor cells already carry the correct shared-text + style in the original sheet.
$donorZero = $ws.Cells.Item(14,3)   # holds text "0" with style s=14
$donorStar = $ws.Cells.Item(14,5)   # holds text "***.*" with style s=14

$donorZero.Copy($ws.Cells.Item(20,3))   # C20
$donorZero.Copy($ws.Cells.Item(22,4))   # D22
$donorStar.Copy($ws.Cells.Item(22,5))   # E22
$donorZero.Copy($ws.Cells.Item(27,6))   # F27

# --- Cells changing from special text markers to numeric values ---
$c = $ws.Cells.Item(15,3)   # C15
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(15,6)   # F15
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(20,4)   # D20
$c.NumberFormat = "#,##0"
$c.Value = 7
$c = $ws.Cells.Item(20,5)   # E20
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100
$c = $ws.Cells.Item(22,6)   # F22
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(26,3)   # C26
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(28,4)   # D28
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(28,5)   # E28
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100
$c = $ws.Cells.Item(29,4)   # D29
$c.NumberFormat = "#,##0"
$c.Value = 1
$c = $ws.Cells.Item(29,5)   # E29
$c.NumberFormat = "#,##0.0;""-""#,##0.0"
$c.Value = -100

# --- Plain numeric value updates ---
$ws.Cells.Item(15,5).Value = -50   # E15
$ws.Cells.Item(15,7).Value = 5   # G15
$ws.Cells.Item(15,8).Value = -80   # H15
$ws.Cells.Item(15,9).Value = 2   # I15
$ws.Cells.Item(15,10).Value = 7   # J15
$ws.Cells.Item(15,11).Value = -71.428571428571   # K15
$ws.Cells.Item(15,12).Value = -33.333333333333   # L15
$ws.Cells.Item(15,13).Value = -33.333333333333   # M15
$ws.Cells.Item(15,14).Value = -88.235294117647   # N15
$ws.Cells.Item(16,3).Value = 1   # C16
$ws.Cells.Item(16,4).Value = 2   # D16
$ws.Cells.Item(16,5).Value = -50   # E16
$ws.Cells.Item(16,6).Value = 17   # F16
$ws.Cells.Item(16,8).Value = 41.666666666666   # H16
$ws.Cells.Item(16,9).Value = 36   # I16
$ws.Cells.Item(16,10).Value = 23   # J16
$ws.Cells.Item(16,11).Value = 56.521739130434   # K16
$ws.Cells.Item(16,12).Value = 80   # L16
$ws.Cells.Item(16,13).Value = -29.411764705882   # M16
$ws.Cells.Item(16,14).Value = -82.524271844660   # N16
$ws.Cells.Item(17,3).Value = 4   # C17
$ws.Cells.Item(17,4).Value = 3   # D17
$ws.Cells.Item(17,5).Value = 33.333333333333   # E17
$ws.Cells.Item(17,6).Value = 14   # F17
$ws.Cells.Item(17,8).Value = -30   # H17
$ws.Cells.Item(17,9).Value = 40   # I17
$ws.Cells.Item(17,10).Value = 48   # J17
$ws.Cells.Item(17,11).Value = -16.666666666666   # K17
$ws.Cells.Item(17,12).Value = -6.976744186046   # L17
$ws.Cells.Item(17,13).Value = -24.528301886792   # M17
$ws.Cells.Item(17,14).Value = -68.75   # N17
$ws.Cells.Item(18,3).Value = 1   # C18
$ws.Cells.Item(18,4).Value = 5   # D18
$ws.Cells.Item(18,9).Value = 32   # I18
$ws.Cells.Item(18,10).Value = 38   # J18
$ws.Cells.Item(18,11).Value = -15.789473684210   # K18
$ws.Cells.Item(18,12).Value = 10.344827586206   # L18
$ws.Cells.Item(18,13).Value = 10.344827586206   # M18
$ws.Cells.Item(18,14).Value = -70.370370370370   # N18
$ws.Cells.Item(19,4).Value = 3   # D19
$ws.Cells.Item(19,5).Value = 66.666666666666   # E19
$ws.Cells.Item(19,6).Value = 18   # F19
$ws.Cells.Item(19,7).Value = 23   # G19
$ws.Cells.Item(19,8).Value = -21.739130434782   # H19
$ws.Cells.Item(19,9).Value = 60   # I19
$ws.Cells.Item(19,10).Value = 58   # J19
$ws.Cells.Item(19,11).Value = 3.448275862068   # K19
$ws.Cells.Item(19,12).Value = 11.111111111111   # L19
$ws.Cells.Item(19,13).Value = 46.341463414634   # M19
$ws.Cells.Item(19,14).Value = 5.263157894736   # N19
$ws.Cells.Item(20,7).Value = 17   # G20
$ws.Cells.Item(20,8).Value = -47.058823529411   # H20
$ws.Cells.Item(20,9).Value = 18   # I20
$ws.Cells.Item(20,10).Value = 34   # J20
$ws.Cells.Item(20,11).Value = -47.058823529411   # K20
$ws.Cells.Item(20,12).Value = 20   # L20
$ws.Cells.Item(20,13).Value = 12.5   # M20
$ws.Cells.Item(20,14).Value = -83.928571428571   # N20
$ws.Cells.Item(21,3).Value = 12   # C21
$ws.Cells.Item(21,4).Value = 22   # D21
$ws.Cells.Item(21,5).Value = -45.454545454545   # E21
$ws.Cells.Item(21,6).Value = 71   # F21
$ws.Cells.Item(21,7).Value = 103   # G21
$ws.Cells.Item(21,8).Value = -31.067961165048   # H21
$ws.Cells.Item(21,9).Value = 188   # I21
$ws.Cells.Item(21,10).Value = 210   # J21
$ws.Cells.Item(21,11).Value = -10.476190476190   # K21
$ws.Cells.Item(21,12).Value = 13.939393939393   # L21
$ws.Cells.Item(21,13).Value = -4.568527918781   # M21
$ws.Cells.Item(21,14).Value = -70.393700787401   # N21
$ws.Cells.Item(22,7).Value = 1   # G22
$ws.Cells.Item(22,8).Value = 0   # H22
$ws.Cells.Item(22,9).Value = 2   # I22
$ws.Cells.Item(22,11).Value = -50   # K22
$ws.Cells.Item(22,12).Value = -33.333333333333   # L22
$ws.Cells.Item(23,6).Value = 4   # F23
$ws.Cells.Item(23,7).Value = 6   # G23
$ws.Cells.Item(23,8).Value = -33.333333333333   # H23
$ws.Cells.Item(23,9).Value = 14   # I23
$ws.Cells.Item(23,10).Value = 13   # J23
$ws.Cells.Item(23,11).Value = 7.692307692307   # K23
$ws.Cells.Item(23,12).Value = -12.5   # L23
$ws.Cells.Item(23,13).Value = -6.666666666666   # M23
$ws.Cells.Item(24,3).Value = 11   # C24
$ws.Cells.Item(24,4).Value = 8   # D24
$ws.Cells.Item(24,5).Value = 37.5   # E24
$ws.Cells.Item(24,6).Value = 67   # F24
$ws.Cells.Item(24,7).Value = 36   # G24
$ws.Cells.Item(24,8).Value = 86.111111111111   # H24
$ws.Cells.Item(24,9).Value = 145   # I24
$ws.Cells.Item(24,10).Value = 102   # J24
$ws.Cells.Item(24,11).Value = 42.156862745098   # K24
$ws.Cells.Item(24,12).Value = 81.25   # L24
$ws.Cells.Item(24,13).Value = 25   # M24
$ws.Cells.Item(25,3).Value = 9   # C25
$ws.Cells.Item(25,4).Value = 6   # D25
$ws.Cells.Item(25,5).Value = 50   # E25
$ws.Cells.Item(25,6).Value = 31   # F25
$ws.Cells.Item(25,7).Value = 26   # G25
$ws.Cells.Item(25,8).Value = 19.230769230769   # H25
$ws.Cells.Item(25,9).Value = 71   # I25
$ws.Cells.Item(25,10).Value = 63   # J25
$ws.Cells.Item(25,11).Value = 12.698412698412   # K25
$ws.Cells.Item(25,12).Value = 61.363636363636   # L25
$ws.Cells.Item(25,13).Value = -41.322314049586   # M25
$ws.Cells.Item(26,5).Value = -50   # E26
$ws.Cells.Item(26,7).Value = 5   # G26
$ws.Cells.Item(26,8).Value = -80   # H26
$ws.Cells.Item(26,9).Value = 4   # I26
$ws.Cells.Item(26,10).Value = 8   # J26
$ws.Cells.Item(26,11).Value = -50   # K26
$ws.Cells.Item(26,12).Value = -20   # L26
$ws.Cells.Item(27,7).Value = 3   # G27
$ws.Cells.Item(27,8).Value = -100   # H27
$ws.Cells.Item(27,10).Value = 7   # J27
$ws.Cells.Item(27,11).Value = -85.714285714285   # K27
$ws.Cells.Item(27,12).Value = -90.909090909090   # L27
$ws.Cells.Item(28,7).Value = 2   # G28
$ws.Cells.Item(28,10).Value = 4   # J28
$ws.Cells.Item(28,11).Value = -50   # K28
$ws.Cells.Item(28,12).Value = 0   # L28
$ws.Cells.Item(28,14).Value = -93.548387096774   # N28
$ws.Cells.Item(29,7).Value = 2   # G29
$ws.Cells.Item(29,10).Value = 4   # J29
$ws.Cells.Item(29,11).Value = -50   # K29
$ws.Cells.Item(29,12).Value = 0   # L29
$ws.Cells.Item(29,14).Value = -93.103448275862   # N29
